# ZBP_07_testovani.xlsx update:
#  - a new survey wave "16. 3. 2021" is added as a new trailing column on both
#    sheets (data: column AA, pocetR: column Z)
#  - the footnote "aktualizace 9. 3. 2021" becomes "aktualizace 23. 3. 2021"
#    on both sheets

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "data" (percentages)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("data")

# New header cell AA1, inheriting the header format (border/bold/alignment)
# used by the rest of row 1 (copy format from Z1, the previous last header).
$ws1.Range("Z1").Copy()
$ws1.Range("AA1").PasteSpecial(-4122)
$ws1.Range("AA1").Value = "16. 3. 2021"

$dataAA = @(0.08,0.08,0.21,0.11,0.08,0.19,0.05,0.1,0.29,0.08,0.05,0.15,0.08,0.1,0.24,0.07000000000000001,0.06,0.19,0.07000000000000001,0.12,0.33,0.13,0.06,0.16,0.07000000000000001,0.02,0.1,0.06,0.05,0.12,0.08,0.1,0.37,0.14,0.11,0.36,0.09,0.08,0.15,0.06,0.22,0.15,0.08,0.11,0.23,0.08,0.08,0.22,0.07000000000000001,0.14,0.41,0.08,0.14,0.23,0.11,0.1,0.21)

for ($i = 0; $i -lt $dataAA.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 27).Value = $dataAA[$i]
}

# Footnote on row 59 (percentages footnote)
$ws1.Range("A59").Value = "Život během pandemie, Testování, % respondentů celkově a ve skupinách, aktualizace 23. 3. 2021"

# ---------------------------------------------------------------------------
# Sheet "pocetR" (sample sizes)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("pocetR")

# New header cell Z1, inheriting the header format used by Y1.
$ws2.Range("Y1").Copy()
$ws2.Range("Z1").PasteSpecial(-4122)
$ws2.Range("Z1").Value = "16. 3. 2021"

$dataZ = @(2101,500,776,825,1015,1086,988,171,640,302,682,105,205,168,386,108,328,180,113)

for ($i = 0; $i -lt $dataZ.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 26).Value = $dataZ[$i]
}

# Row 21 is a blank spacer row (all cells hold an empty string); extend it to
# column Z to match the rest of the row.
$ws2.Range("Y21").Copy()
$ws2.Range("Z21").PasteSpecial(-4163)

# Footnote on row 21 (sample-size footnote)
$ws2.Range("A21").Value = "Život během pandemie, Testování, velikost dotázaného souboru celkově a ve skupinách, aktualizace 23. 3. 2021"
